$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out the extent of the existing data (header row + player rows).
$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1
$lastDataCol = $used.Column + $used.Columns.Count - 1

# New columns go right after the last existing column (AC -> AD, AE, AF).
$winsCol = $lastDataCol + 1
$lossesCol = $lastDataCol + 2
$tiesCol = $lastDataCol + 3

# Header row: add Wins / Losses / Ties headers, matching the existing header style.
$ws.Cells.Item(1, $winsCol).Value = "Wins"
$ws.Cells.Item(1, $lossesCol).Value = "Losses"
$ws.Cells.Item(1, $tiesCol).Value = "Ties"

$headerRange = $ws.Range($ws.Cells.Item(1, $winsCol), $ws.Cells.Item(1, $tiesCol))
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Data rows: team's win/loss/tie record, repeated for every player row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $winsCol).Value = 73
    $ws.Cells.Item($r, $lossesCol).Value = 89
    $ws.Cells.Item($r, $tiesCol).Value = 0
}
